$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above row 121 (existing rows 121-189 shift down to 123-191).
$ws.Rows.Item(121).Resize(2).Insert()

# --- New row 121: Red Globe / Provincia del Elqui / $/bandeja 12 kilos ---
$ws.Cells.Item(121, 1).Value  = 4
$ws.Cells.Item(121, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(121, 3).Value  = "Los Lagos"
$ws.Cells.Item(121, 4).Value  = 44582
$ws.Cells.Item(121, 5).Value  = 10
$ws.Cells.Item(121, 6).Value  = "Fruta"
$ws.Cells.Item(121, 7).Value  = 100109
$ws.Cells.Item(121, 8).Value  = "Uva"
$ws.Cells.Item(121, 9).Value  = 100109001
$ws.Cells.Item(121, 10).Value = "Uva"
$ws.Cells.Item(121, 11).Value = "Red Globe"
$ws.Cells.Item(121, 12).Value = "Primera"
$ws.Cells.Item(121, 13).Value = 300
$ws.Cells.Item(121, 14).Value = 15000
$ws.Cells.Item(121, 15).Value = 16000
$ws.Cells.Item(121, 16).Value = 15500
$ws.Cells.Item(121, 17).Value = "$/bandeja 12 kilos"
$ws.Cells.Item(121, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(121, 19).Value = 1292
$ws.Cells.Item(121, 20).Value = 12

# --- New row 122: Superior Seedless / Provincia de Limarí / $/caja 15 kilos ---
$ws.Cells.Item(122, 1).Value  = 4
$ws.Cells.Item(122, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(122, 3).Value  = "Los Lagos"
$ws.Cells.Item(122, 4).Value  = 44582
$ws.Cells.Item(122, 5).Value  = 10
$ws.Cells.Item(122, 6).Value  = "Fruta"
$ws.Cells.Item(122, 7).Value  = 100109
$ws.Cells.Item(122, 8).Value  = "Uva"
$ws.Cells.Item(122, 9).Value  = 100109001
$ws.Cells.Item(122, 10).Value = "Uva"
$ws.Cells.Item(122, 11).Value = "Superior Seedless"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 400
$ws.Cells.Item(122, 14).Value = 19000
$ws.Cells.Item(122, 15).Value = 20000
$ws.Cells.Item(122, 16).Value = 19500
$ws.Cells.Item(122, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(122, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(122, 19).Value = 1300
$ws.Cells.Item(122, 20).Value = 15
